$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data (refresh scrape).
# Cells whose new text looks like a plain number are forced to "Text"
# number format first, so Excel stores the exact original string
# (preserving trailing zeros / tiny decimals) instead of silently
# re-parsing it into a numeric value.

$ws.Range('D2').Value = '68.865.34'
$ws.Range('E2').Value = '  +4.59%  '
$ws.Range('D3').Value = '3.378.79'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.86'
$ws.Range('E5').Value = '  +6.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.59'
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.599'
$ws.Range('E7').Value = '  +4.11%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +5.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.587'
$ws.Range('E10').Value = '  +1.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '47.36'
$ws.Range('E11').Value = '  +3.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000281'
$ws.Range('E12').Value = '  +7.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '640.76'
$ws.Range('E13').Value = '  +12.42%  '
$ws.Range('D14').Value = '3.915.57'
$ws.Range('E14').Value = '  +1.73%  '
$ws.Range('E15').Value = '  +1.36%  '
$ws.Range('D16').Value = '69.010.03'
$ws.Range('E16').Value = '  +4.83%  '
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('D18').Value = '3.370.92'
$ws.Range('E18').Value = '  +1.53%  '
$ws.Range('E19').Value = '  +1.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.11'
$ws.Range('E20').Value = '  +2.73%  '
$ws.Range('E21').Value = '  +2.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.95'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.09'
$ws.Range('E23').Value = '  +1.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '99.26'
$ws.Range('E24').Value = '  +1.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.10'
$ws.Range('E25').Value = '  +4.27%  '
$ws.Range('E26').Value = '  +5.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.81'
$ws.Range('E27').Value = '  +4.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '33.00'
$ws.Range('E28').Value = '  +8.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.70'
$ws.Range('E29').Value = '  +2.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.83'
$ws.Range('E30').Value = '  +2.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '614.47'
$ws.Range('E31').Value = '  +9.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.70'
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').Value = '3.991.64'
$ws.Range('E33').Value = '  +6.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.11'
$ws.Range('E34').Value = '  +2.71%  '
$ws.Range('E35').Value = '  +2.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.11'
$ws.Range('E37').Value = '  +1.07%  '
$ws.Range('E38').Value = '  +8.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.32'
$ws.Range('E39').Value = '  +6.78%  '
$ws.Range('E40').Value = '  +4.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '33.71'
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('D42').Value = '0.0₃0708'
$ws.Range('E42').Value = '  +3.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.42'
$ws.Range('E43').Value = '  +2.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.344'
$ws.Range('E44').Value = '  +3.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0424'
$ws.Range('E45').Value = '  +3.95%  '
$ws.Range('E46').Value = '  +2.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.59'
$ws.Range('E47').Value = '  +3.42%  '
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('E49').Value = '  +9.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '131.02'
$ws.Range('E50').Value = '  +5.07%  '
$ws.Range('E51').Value = '  +7.44%  '
